# Update room names in the "Evan QI_CampA_timetable" worksheet:
#   "Room Stephane" -> "Room G19"
#   "Room Ivy"      -> "Room G14"
#   "Room Shelley"  -> "Room G22"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value  = "Private Lesson with Stephane RETY `n(Room G19)"
$ws.Range("E7").Value  = "Private Lesson with Stephane RETY & pianist `n(Room G19)"
$ws.Range("F7").Value  = "Flute MasterClass`n(Room G19)"

$ws.Range("D11").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"

$ws.Range("C19").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("D19").Value = "Rehearsal with pianist`n(Room G22)"
$ws.Range("F19").Value = "Flute MasterClass`n(Room G19)"

$ws.Range("B27").Value = "Ensemble `n(Room G14)"
$ws.Range("C27").Value = "Ensemble `n(Room G14)"
$ws.Range("D27").Value = "Ensemble `n(Room G14)"
$ws.Range("E27").Value = "Ensemble `n(Room G14)"
$ws.Range("F27").Value = "Ensemble `n(Room G14)"
